# Add the new "TMF8801" worksheet (I2C Time-of-Flight sensor register map) as
# the last tab, matching the commit: "added TMF8801 code to demonstrate
# isolated i2c".

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the current last sheet (AD7961) so it lands at
# the end of the tab strip.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "TMF8801"

# Header row - same column headers used by every other register sheet in
# this workbook.
$ws.Cells.Item(1, 1).Value = "Name"
$ws.Cells.Item(1, 2).Value = "Hex Address"
$ws.Cells.Item(1, 3).Value = "Default Value"
$ws.Cells.Item(1, 4).Value = "Bit Width"
$ws.Cells.Item(1, 5).Value = "Bit Index (High)"
$ws.Cells.Item(1, 6).Value = "Bit Index (Low)"

# TMF8801 register map: Name, Hex Address, Default Value, Bit Width,
# Bit Index (High), Bit Index (Low)
$data = @(
    ,@("APPID","0x00","0xff",16,15,0)
    ,@("APPREV_MAJOR","0x01","0x00",16,15,0)
    ,@("APPREQID","0x02","0xff",16,15,0)
    ,@("ENABLE","0xE0","0x00",16,15,0)
    ,@("INT_STATUS","0xE1","0x00",16,15,0)
    ,@("INT_ENAB","0xE2","0x00",16,15,0)
    ,@("ID","0xE3","0x00",16,15,0)
    ,@("REVID","0xE4","0x00",16,15,0)
    ,@("CMD_DATA9","0x06","0x00",16,15,0)
    ,@("CMD_DATA8","0x07","0x00",16,15,0)
    ,@("CMD_DATA7","0x08","0x00",16,15,0)
    ,@("CMD_DATA6","0x09","0x00",16,15,0)
    ,@("CMD_DATA5","0x0A","0x00",16,15,0)
    ,@("CMD_DATA4","0x0B","0x00",16,15,0)
    ,@("CMD_DATA3","0x0C","0x00",16,15,0)
    ,@("CMD_DATA2","0x0D","0x00",16,15,0)
    ,@("CMD_DATA1","0x0E","0x00",16,15,0)
    ,@("CMD_DATA0","0x0F","0x00",16,15,0)
    ,@("COMMAND","0x10","0x00",16,15,0)
    ,@("PREVIOUS","0x11","0x00",16,15,0)
    ,@("APPREV_MINOR","0x12","0x00",16,15,0)
    ,@("APPREV_PATCH","0x13","0x00",16,15,0)
    ,@("STATUS","0x1D","0x00",16,15,0)
    ,@("REGISTER_CONTENTS","0x1E","0x00",16,15,0)
    ,@("TID","0x1F","0x00",16,15,0)
    ,@("HIST_START","0x20","0x00",16,15,0)
    ,@("HIST_END","0x9F","0x00",16,15,0)
    ,@("SERIES_NUMBER_0","0x28","0x00",16,15,0)
    ,@("SERIES_NUMBER_0","0x29","0x00",16,15,0)
    ,@("IDENTIFICATION_NUMBER_0","0x2A","0x00",16,15,0)
    ,@("IDENTIFICATION_NUMBER_1","0x2B","0x00",16,15,0)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# Column widths roughly matching the authored sheet.
$ws.Columns.Item(1).ColumnWidth = 31.83203125
$ws.Columns.Item(5).ColumnWidth = 14.6640625

# Leave the cursor where the author left it.
$ws.Range("F23").Select()
